$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2724.9866
$ws.Range("I15").Value = 2724.9866
$ws.Range("K15").Value = 8174.959800000001
$ws.Range("M15").Value = -8005.959800000001
$ws.Range("H32").Value = 1142.2142
$ws.Range("I32").Value = 533.3333
$ws.Range("J32").Value = 1308.2727
$ws.Range("K32").Value = 533.3333
$ws.Range("L32").Value = 1308.2727
$ws.Range("M32").Value = -207.3333
$ws.Range("N32").Value = -1960.2727
$ws.Range("H58").Value = 1003.0625
$ws.Range("I58").Value = 475.63635
$ws.Range("J58").Value = 2163.4
$ws.Range("K58").Value = 1426.90905
$ws.Range("L58").Value = 6490.200000000001
$ws.Range("M58").Value = -1276.90905
$ws.Range("N58").Value = -6790.200000000001
$ws.Range("H64").Value = 4299.3335
$ws.Range("I64").Value = 4249.1665
$ws.Range("J64").Value = 4500
$ws.Range("K64").Value = 4249.1665
$ws.Range("L64").Value = 4500
$ws.Range("M64").Value = -4001.1665
$ws.Range("N64").Value = -4996
$ws.Range("H67").Value = 4299.3335
$ws.Range("I67").Value = 4249.1665
$ws.Range("J67").Value = 4500
$ws.Range("K67").Value = 4249.1665
$ws.Range("L67").Value = 4500
$ws.Range("M67").Value = -3391.1665
$ws.Range("N67").Value = -6216
$ws.Range("H76").Value = 3816.2778
$ws.Range("I76").Value = 4308.4546
$ws.Range("J76").Value = 3042.8572
$ws.Range("K76").Value = 4308.4546
$ws.Range("L76").Value = 3042.8572
$ws.Range("M76").Value = -3993.4546
$ws.Range("N76").Value = -3672.8572
$ws.Range("H79").Value = 3816.2778
$ws.Range("I79").Value = 4308.4546
$ws.Range("J79").Value = 3042.8572
$ws.Range("K79").Value = 4308.4546
$ws.Range("L79").Value = 3042.8572
$ws.Range("M79").Value = -3216.4546
$ws.Range("N79").Value = -5226.8572
$ws.Range("H92").Value = 1871.8182
$ws.Range("J92").Value = 1625
$ws.Range("L92").Value = 1625
$ws.Range("N92").Value = -4121
$ws.Range("H116").Value = 2698.2666
$ws.Range("I116").Value = 1718.8889
$ws.Range("K116").Value = 1718.8889
$ws.Range("M116").Value = 1723.1111
$ws.Range("H132").Value = 6294655.5
$ws.Range("I132").Value = 10106434
$ws.Range("J132").Value = 5219.8
$ws.Range("K132").Value = 30319302
$ws.Range("L132").Value = 15659.4
$ws.Range("M132").Value = -30316772
$ws.Range("N132").Value = -20719.4
$ws.Range("H135").Value = 21739586
$ws.Range("I135").Value = 208.925
$ws.Range("J135").Value = 166668770
$ws.Range("K135").Value = 1880.325
$ws.Range("L135").Value = 1500018930
$ws.Range("M135").Value = 654.675
$ws.Range("N135").Value = -1500024000
$ws.Range("H137").Value = 1200.9166
$ws.Range("I137").Value = 877.3333
$ws.Range("K137").Value = 2631.9999
$ws.Range("M137").Value = -81.9998999999998
$ws.Range("H138").Value = 1285.61
$ws.Range("I138").Value = 636.5
$ws.Range("J138").Value = 1591.0735
$ws.Range("K138").Value = 1909.5
$ws.Range("L138").Value = 4773.220499999999
$ws.Range("M138").Value = 3230.5
$ws.Range("N138").Value = -15053.2205

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4337.9404
$ws.Range("I32").Value = 3819.6528
$ws.Range("J32").Value = 7447.6665
$ws.Range("K32").Value = 3819.6528
$ws.Range("L32").Value = 7447.6665
$ws.Range("M32").Value = -3532.6528
$ws.Range("N32").Value = -8021.6665
$ws.Range("H45").Value = 1501.5555
$ws.Range("I45").Value = 1666.6666
$ws.Range("J45").Value = 1171.3334
$ws.Range("K45").Value = 1666.6666
$ws.Range("L45").Value = 1171.3334
$ws.Range("M45").Value = -1289.6666
$ws.Range("N45").Value = -1925.3334
$ws.Range("H61").Value = 25642172
$ws.Range("I61").Value = 32258830
$ws.Range("J61").Value = 2625
$ws.Range("K61").Value = 32258830
$ws.Range("L61").Value = 2625
$ws.Range("M61").Value = -32258618
$ws.Range("N61").Value = -3049
$ws.Range("H74").Value = 972.0208
$ws.Range("I74").Value = 753.8372000000001
$ws.Range("J74").Value = 2848.4
$ws.Range("K74").Value = 753.8372000000001
$ws.Range("L74").Value = 2848.4
$ws.Range("M74").Value = 120.1627999999999
$ws.Range("N74").Value = -4596.4
$ws.Range("H77").Value = 972.0208
$ws.Range("I77").Value = 753.8372000000001
$ws.Range("J77").Value = 2848.4
$ws.Range("K77").Value = 3769.186
$ws.Range("L77").Value = 14242
$ws.Range("M77").Value = 598.8139999999999
$ws.Range("N77").Value = -22978
$ws.Range("H92").Value = 1903333.4
$ws.Range("J92").Value = 1903333.4
$ws.Range("L92").Value = 1903333.4
$ws.Range("N92").Value = -1908325.4
$ws.Range("H97").Value = 275.47058
$ws.Range("I97").Value = 298.77777
$ws.Range("J97").Value = 185.57143
$ws.Range("K97").Value = 298.77777
$ws.Range("L97").Value = 185.57143
$ws.Range("M97").Value = 197.22223
$ws.Range("N97").Value = -1177.57143
$ws.Range("H102").Value = 11113128
$ws.Range("I102").Value = 13890673
$ws.Range("J102").Value = 2947.3333
$ws.Range("K102").Value = 13890673
$ws.Range("L102").Value = 2947.3333
$ws.Range("M102").Value = -13889051
$ws.Range("N102").Value = -6191.3333
$ws.Range("H110").Value = 1183.9354
$ws.Range("I110").Value = 746.8823
$ws.Range("K110").Value = 746.8823
$ws.Range("M110").Value = 1298.1177
$ws.Range("H132").Value = 2211.1538
$ws.Range("I132").Value = 1986.2
$ws.Range("K132").Value = 5958.6
$ws.Range("M132").Value = -3428.6
$ws.Range("H136").Value = 25642172
$ws.Range("I136").Value = 32258830
$ws.Range("J136").Value = 2625
$ws.Range("K136").Value = 96776490
$ws.Range("L136").Value = 7875
$ws.Range("M136").Value = -96773940
$ws.Range("N136").Value = -12975

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 15625798
$ws.Range("I94").Value = 35714544
$ws.Range("J94").Value = 1218.4445
$ws.Range("K94").Value = 35714544
$ws.Range("L94").Value = 1218.4445
$ws.Range("M94").Value = -35714093
$ws.Range("N94").Value = -2120.4445

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2657.0557
$ws.Range("I31").Value = 2636.8823
$ws.Range("K31").Value = 2636.8823
$ws.Range("M31").Value = -2341.8823
$ws.Range("H34").Value = 2657.0557
$ws.Range("I34").Value = 2636.8823
$ws.Range("K34").Value = 2636.8823
$ws.Range("M34").Value = -2434.8823
$ws.Range("H58").Value = 894.84784
$ws.Range("I58").Value = 820.875
$ws.Range("K58").Value = 820.875
$ws.Range("M58").Value = -617.875
$ws.Range("H86").Value = 14709.27
$ws.Range("I86").Value = 8858
$ws.Range("J86").Value = 24071.3
$ws.Range("K86").Value = 8858
$ws.Range("L86").Value = 24071.3
$ws.Range("M86").Value = -7735
$ws.Range("N86").Value = -26317.3
$ws.Range("H89").Value = 14709.27
$ws.Range("I89").Value = 8858
$ws.Range("J89").Value = 24071.3
$ws.Range("K89").Value = 44290
$ws.Range("L89").Value = 120356.5
$ws.Range("M89").Value = -38674
$ws.Range("N89").Value = -131588.5
$ws.Range("H109").Value = 14544.556
$ws.Range("J109").Value = 14544.556
$ws.Range("L109").Value = 14544.556
$ws.Range("N109").Value = -16624.556
$ws.Range("H132").Value = 2538.682
$ws.Range("I132").Value = 2202.7856
$ws.Range("J132").Value = 3126.5
$ws.Range("K132").Value = 6608.3568
$ws.Range("L132").Value = 9379.5
$ws.Range("M132").Value = -4078.3568
$ws.Range("N132").Value = -14439.5
$ws.Range("H134").Value = 16667876
$ws.Range("I134").Value = 1149.25
$ws.Range("J134").Value = 50001330
$ws.Range("K134").Value = 3447.75
$ws.Range("L134").Value = 150003990
$ws.Range("M134").Value = -912.75
$ws.Range("N134").Value = -150009060
$ws.Range("H136").Value = 894.84784
$ws.Range("I136").Value = 820.875
$ws.Range("K136").Value = 2462.625
$ws.Range("M136").Value = 87.375

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2977.762
$ws.Range("J39").Value = 2840.7222
$ws.Range("L39").Value = 8522.1666
$ws.Range("N39").Value = -9110.1666
$ws.Range("H55").Value = 2415
$ws.Range("J55").Value = 3181
$ws.Range("L55").Value = 9543
$ws.Range("N55").Value = -9897

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1803
$ws.Range("I113").Value = 1784.2
$ws.Range("J113").Value = 1850
$ws.Range("K113").Value = 1784.2
$ws.Range("L113").Value = 1850
$ws.Range("M113").Value = 385.8
$ws.Range("N113").Value = -6190
$ws.Range("H126").Value = 2966.6667
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 2966.6667
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 8900.000100000001
$ws.Range("N126").Value = -13840.0001
$ws.Range("H132").Value = 2318.5454
$ws.Range("I132").Value = 2158.4736
$ws.Range("J132").Value = 2535.7856
$ws.Range("K132").Value = 6475.4208
$ws.Range("L132").Value = 7607.3568
$ws.Range("M132").Value = -3945.4208
$ws.Range("N132").Value = -12667.3568

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2023.4117
$ws.Range("I82").Value = 2049.9167
$ws.Range("J82").Value = 1959.8
$ws.Range("K82").Value = 2049.9167
$ws.Range("L82").Value = 1959.8
$ws.Range("M82").Value = -1688.9167
$ws.Range("N82").Value = -2681.8
$ws.Range("H85").Value = 2023.4117
$ws.Range("I85").Value = 2049.9167
$ws.Range("J85").Value = 1959.8
$ws.Range("K85").Value = 2049.9167
$ws.Range("L85").Value = 1959.8
$ws.Range("M85").Value = -801.9167000000002
$ws.Range("N85").Value = -4455.8
$ws.Range("H93").Value = 613.125
$ws.Range("I93").Value = 543.5714
$ws.Range("J93").Value = 1100
$ws.Range("K93").Value = 543.5714
$ws.Range("L93").Value = 1100
$ws.Range("M93").Value = 704.4286
$ws.Range("N93").Value = -3596
$ws.Range("H94").Value = 49999.5
$ws.Range("J94").Value = 49999.5
$ws.Range("L94").Value = 49999.5
$ws.Range("N94").Value = -51351.5
$ws.Range("H123").Value = 40930
$ws.Range("J123").Value = 40930
$ws.Range("L123").Value = 40930
$ws.Range("N123").Value = -50730
$ws.Range("H132").Value = 20796.191
$ws.Range("I132").Value = 940.0333000000001
$ws.Range("K132").Value = 2820.0999
$ws.Range("M132").Value = -290.0999000000002
$ws.Range("H136").Value = 1972.1875
$ws.Range("I136").Value = 1862.9166
$ws.Range("J136").Value = 2300
$ws.Range("K136").Value = 5588.7498
$ws.Range("L136").Value = 6900
$ws.Range("M136").Value = -3038.7498
$ws.Range("N136").Value = -12000

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1134.2084
$ws.Range("I132").Value = 712.14
$ws.Range("J132").Value = 2093.4546
$ws.Range("K132").Value = 2136.42
$ws.Range("L132").Value = 6280.3638
$ws.Range("M132").Value = 393.5799999999999
$ws.Range("N132").Value = -11340.3638

# ---- Special case: GSM M126 cell fully removed ----
$wsGSM = $wb.Worksheets.Item("GSM")
$wsGSM.Range("M126").ClearContents()
